$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 986.36365
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 1023.8095
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 3071.4285
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -4063.4285

$ws.Range("H98").Value = 936.46155
$ws.Range("I98").Value = 743.1818
$ws.Range("K98").Value = 743.1818
$ws.Range("M98").Value = 754.8182

$ws.Range("H122").Value = 936.46155
$ws.Range("I122").Value = 743.1818
$ws.Range("K122").Value = 2229.5454
$ws.Range("M122").Value = 220.4546

$ws.Range("H129").Value = 173381.69
$ws.Range("I129").Value = 324.25
$ws.Range("J129").Value = 186200.77
$ws.Range("K129").Value = 972.75
$ws.Range("L129").Value = 558602.3099999999
$ws.Range("M129").Value = 4027.25
$ws.Range("N129").Value = -568602.3099999999

$ws.Range("H132").Value = 2785.6765
$ws.Range("I132").Value = 2824.4546
$ws.Range("K132").Value = 8473.363799999999
$ws.Range("M132").Value = -5943.363799999999

$ws.Range("H138").Value = 2407.1143
$ws.Range("I138").Value = 1526.9412
$ws.Range("K138").Value = 4580.8236
$ws.Range("M138").Value = 559.1764000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4287.3013
$ws.Range("I32").Value = 3664.7883
$ws.Range("K32").Value = 3664.7883
$ws.Range("M32").Value = -3377.7883

$ws.Range("H132").Value = 10413.018
$ws.Range("I132").Value = 1559.2291
$ws.Range("J132").Value = 52911.2
$ws.Range("K132").Value = 4677.6873
$ws.Range("L132").Value = 158733.6
$ws.Range("M132").Value = -2147.6873
$ws.Range("N132").Value = -163793.6

$ws.Range("H135").Value = 35881.5
$ws.Range("J135").Value = 35881.5
$ws.Range("L135").Value = 35881.5
$ws.Range("N135").Value = -46021.5

$ws.Range("H139").Value = 40594.25
$ws.Range("J139").Value = 40594.25
$ws.Range("L139").Value = 40594.25
$ws.Range("N139").Value = -50874.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 9482.9
$ws.Range("I54").Value = 2365.8
$ws.Range("J54").Value = 16600
$ws.Range("K54").Value = 2365.8
$ws.Range("L54").Value = 16600
$ws.Range("M54").Value = -1881.8
$ws.Range("N54").Value = -17568

$ws.Range("H70").Value = 93050
$ws.Range("J70").Value = 93050
$ws.Range("L70").Value = 93050
$ws.Range("N70").Value = -93636

$ws.Range("H73").Value = 93050
$ws.Range("J73").Value = 93050
$ws.Range("L73").Value = 93050
$ws.Range("N73").Value = -95078

$ws.Range("H105").Value = 1831.4348
$ws.Range("I105").Value = 1595.091
$ws.Range("K105").Value = 1595.091
$ws.Range("M105").Value = 151.9090000000001

$ws.Range("H108").Value = 42785
$ws.Range("J108").Value = 42785
$ws.Range("L108").Value = 42785
$ws.Range("N108").Value = -50465

$ws.Range("H134").Value = 3197.718
$ws.Range("I134").Value = 3293.1516
$ws.Range("J134").Value = 2672.8333
$ws.Range("K134").Value = 9879.4548
$ws.Range("L134").Value = 8018.499899999999
$ws.Range("M134").Value = -7344.4548
$ws.Range("N134").Value = -13088.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1301.6818
$ws.Range("J107").Value = 2009.5454
$ws.Range("L107").Value = 2009.5454
$ws.Range("N107").Value = -5849.5454

$ws.Range("H120").Value = 11333.25
$ws.Range("I120").Value = 9888.666999999999
$ws.Range("J120").Value = 12200
$ws.Range("K120").Value = 9888.666999999999
$ws.Range("L120").Value = 12200
$ws.Range("M120").Value = -6259.666999999999
$ws.Range("N120").Value = -19458

$ws.Range("H121").Value = 7793.8335
$ws.Range("I121").Value = 5340
$ws.Range("J121").Value = 20063
$ws.Range("K121").Value = 5340
$ws.Range("L121").Value = 20063
$ws.Range("M121").Value = -4030
$ws.Range("N121").Value = -22683

$ws.Range("H122").Value = 1065
$ws.Range("I122").Value = 1043.75
$ws.Range("J122").Value = 1150
$ws.Range("K122").Value = 3131.25
$ws.Range("L122").Value = 3450
$ws.Range("M122").Value = -681.25
$ws.Range("N122").Value = -8350

$ws.Range("H124").Value = 5405.5
$ws.Range("I124").Value = 2648
$ws.Range("J124").Value = 8163
$ws.Range("K124").Value = 2648
$ws.Range("L124").Value = 8163
$ws.Range("M124").Value = -193
$ws.Range("N124").Value = -13073

$ws.Range("H132").Value = 2055.2368
$ws.Range("I132").Value = 1412.3214
$ws.Range("J132").Value = 3855.4
$ws.Range("K132").Value = 4236.9642
$ws.Range("L132").Value = 11566.2
$ws.Range("M132").Value = -1706.9642
$ws.Range("N132").Value = -16626.2

$ws.Range("H133").Value = 34779.855
$ws.Range("J133").Value = 34779.855
$ws.Range("L133").Value = 34779.855
$ws.Range("N133").Value = -39839.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 717.8
$ws.Range("J122").Value = 884.7
$ws.Range("L122").Value = 7962.3
$ws.Range("N122").Value = -12862.3

$ws.Range("H131").Value = 700.15
$ws.Range("J131").Value = 700.15
$ws.Range("L131").Value = 2100.45
$ws.Range("N131").Value = -12180.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 74076600
$ws.Range("I122").Value = 25642366
$ws.Range("J122").Value = 200005600
$ws.Range("K122").Value = 76927098
$ws.Range("L122").Value = 600016800
$ws.Range("M122").Value = -76924648
$ws.Range("N122").Value = -600021700

$ws.Range("H126").Value = 2947.9546
$ws.Range("I126").Value = 1897.4
$ws.Range("J126").Value = 5199.143
$ws.Range("K126").Value = 5692.200000000001
$ws.Range("L126").Value = 15597.429
$ws.Range("M126").Value = -3222.200000000001
$ws.Range("N126").Value = -20537.429

$ws.Range("H132").Value = 13691.2
$ws.Range("I132").Value = 2490.3901
$ws.Range("K132").Value = 7471.1703
$ws.Range("M132").Value = -4941.1703

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3212.4666
$ws.Range("I22").Value = 5743
$ws.Range("J22").Value = 998.25
$ws.Range("K22").Value = 5743
$ws.Range("L22").Value = 998.25
$ws.Range("M22").Value = -5448
$ws.Range("N22").Value = -1588.25

$ws.Range("H27").Value = 3212.4666
$ws.Range("I27").Value = 5743
$ws.Range("J27").Value = 998.25
$ws.Range("K27").Value = 5743
$ws.Range("L27").Value = 998.25
$ws.Range("M27").Value = -5636
$ws.Range("N27").Value = -1212.25

$ws.Range("H40").Value = 4128
$ws.Range("I40").Value = 3119
$ws.Range("K40").Value = 3119
$ws.Range("M40").Value = -2983

$ws.Range("H93").Value = 1305.75
$ws.Range("I93").Value = 1198.6666
$ws.Range("J93").Value = 1627
$ws.Range("K93").Value = 1198.6666
$ws.Range("L93").Value = 1627
$ws.Range("M93").Value = 49.33339999999998
$ws.Range("N93").Value = -4123

$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("L109").Value = 27000
$ws.Range("N109").Value = -29774

$ws.Range("H132").Value = 432499.06
$ws.Range("I132").Value = 635557.7
$ws.Range("J132").Value = 3819.7778
$ws.Range("K132").Value = 1906673.1
$ws.Range("L132").Value = 11459.3334
$ws.Range("M132").Value = -1904143.1
$ws.Range("N132").Value = -16519.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 14701
$ws.Range("J101").Value = 14701
$ws.Range("L101").Value = 14701
$ws.Range("N101").Value = -21191

$ws.Range("H122").Value = 1299.8572
$ws.Range("I122").Value = 1209.1428
$ws.Range("J122").Value = 1572
$ws.Range("K122").Value = 3627.4284
$ws.Range("L122").Value = 4716
$ws.Range("M122").Value = -1177.4284
$ws.Range("N122").Value = -9616

$ws.Range("H132").Value = 1555.5555
$ws.Range("I132").Value = 1153.1177
$ws.Range("J132").Value = 2239.7
$ws.Range("K132").Value = 3459.3531
$ws.Range("L132").Value = 6719.099999999999
$ws.Range("M132").Value = -929.3531000000003
$ws.Range("N132").Value = -11779.1
